$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data table's last row (43) holds 2025-10-21 (serial 45951); extend
# it with two more rows for the next day, 2025-10-22 (serial 45952),
# one per station, matching the existing row-by-row layout.

# Carry the formatting (date/number styles) of the last existing row down
# into the two new rows before writing their values.
$ws.Range("A43:F43").Copy()
$ws.Range("A44:F45").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(44, 1).Value = 45952
$ws.Cells.Item(44, 2).Value = "四方坪站"
$ws.Cells.Item(44, 3).Value = 8977.26
$ws.Cells.Item(44, 4).Value = 7303.42
$ws.Cells.Item(44, 5).Value = 3108.35
$ws.Cells.Item(44, 6).Value = 388

$ws.Cells.Item(45, 1).Value = 45952
$ws.Cells.Item(45, 2).Value = "高岭站"
$ws.Cells.Item(45, 3).Value = 5175.03
$ws.Cells.Item(45, 4).Value = 4382.98
$ws.Cells.Item(45, 5).Value = 1404.57
$ws.Cells.Item(45, 6).Value = 197

# The sheet had scrolled down (topLeftCell A31) with I38 selected; restore
# the view to the top and select J11.
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("J11").Select()
